$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiple purchase price / quantity-per-unit values stored as text (comma separated)
$ws.Range("H2").Value = "2000 , 200 , 1000"
$ws.Range("J2").Value = "15 , 10 , 10"

# Column width adjustments to fit the new text values
$ws.Columns.Item(8).ColumnWidth = 23.166666666666668
$ws.Columns.Item(10).ColumnWidth = 12.5

# Update the active selection
$ws.Range("J8").Select()
